{"js": "// Recolor the \"AngularJS Storage\" (SNO 14) and \"NodeJS\" (SNO 16) syllabus\n// rows from black to green (#00B050), matching the rest of the topics that\n// are already covered/highlighted in the table.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\n\n// Load each row's cell count and first-cell text so we can locate the two\n// target rows (SNO \"14\" and SNO \"16\") without depending on a fixed index.\nrows.forEach((row) => row.load(\"cellCount,values\"));\nawait context.sync();\n\nconst GREEN = \"#00B050\";\nconst targetSnos = new Set([\"14\", \"16\"]);\n\nfor (let i = 0; i < rows.length; i++) {\n  const row = rows[i];\n  const firstCellText = row.values && row.values[0] ? row.values[0][0].trim() : \"\";\n  if (!targetSnos.has(firstCellText)) {\n    continue;\n  }\n\n  // Recolor every real cell in this row (vertically-merged continuation\n  // cells are not exposed as part of cellCount, so this naturally skips\n  // them - matching the diff, which leaves those untouched).\n  for (let c = 0; c < row.cellCount; c++) {\n    const cell = table.getCell(i, c);\n    cell.body.font.color = GREEN;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Recolor the \"AngularJS Storage\" (SNO 14) and \"NodeJS\" (SNO 16) syllabus\n# rows from black to green (00B050), matching the rest of the topics that\n# are already covered/highlighted in the table.\n#\n# wdColor values are packed as 0x00BBGGRR, so RGB(00,B0,50) -> 0x0050B000.\n$GREEN = 0x0050B000\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$targetSnos = @(\"14\", \"16\")\n\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    $row = $table.Rows.Item($r)\n    $firstCell = $row.Cells.Item(1)\n    $sno = $firstCell.Range.Text.TrimEnd([char]13, [char]7).Trim()\n\n    if ($targetSnos -notcontains $sno) {\n        continue\n    }\n\n    # Recolor every real cell in this row. Vertically-merged continuation\n    # cells are not part of Cells.Count for this row, so this naturally\n    # skips them - matching the diff, which leaves those untouched.\n    for ($c = 1; $c -le $row.Cells.Count; $c++) {\n        $cell = $row.Cells.Item($c)\n        $cell.Range.Font.Color = $GREEN\n    }\n}\n"}
